{"js": "// Update the French \"Th\u00e9 Chai\" slogans to the revised \"Th\u00e9 cha\u00ef\" copy.\n// Each slogan is its own single-run list paragraph, so a straight\n// text-for-text swap per paragraph is the safest approach (keeps the\n// run's formatting/rPr intact while matching the exact new wording).\nconst replacements = [\n  [\n    \"Th\u00e9 Chai\\u00A0: L\\u2019\u00e9pices de la vie\",\n    \"Th\u00e9 cha\u00ef\\u00A0: L\\u2019\u00e9pice de la vie\",\n  ],\n  [\n    \"Th\u00e9 Chai\\u00A0: Un monde de saveur dans une tasse\",\n    \"Th\u00e9 cha\u00ef\\u00A0: Un monde de saveur dans une tasse\",\n  ],\n  [\n    \"Th\u00e9 Chai\\u00A0: D\u00e9couvrir la magie de l\\u2019Inde\",\n    \"Th\u00e9 cha\u00ef\\u00A0: D\u00e9couvrez la magie de l\\u2019Inde\",\n  ],\n  [\n    \"Th\u00e9 Chai\\u00A0: Le m\u00e9lange parfait de sant\u00e9 et de plaisir\",\n    \"Th\u00e9 cha\u00ef\\u00A0: L\\u2019\u00e9quilibre parfait du bien-\u00eatre et du plaisir\",\n  ],\n  [\n    \"Th\u00e9 Chai\\u00A0: Plus que le th\u00e9, un mode de vie\",\n    \"Th\u00e9 cha\u00ef\\u00A0: Plus qu\\u2019un th\u00e9, un mode de vie\",\n  ],\n  [\n    \"Th\u00e9 Chai\\u00A0: Une boisson pour toutes les saisons et les raisons\",\n    \"Th\u00e9 cha\u00ef\\u00A0: Une boisson pour toutes les saisons et toutes les raisons\",\n  ],\n  [\n    \"Th\u00e9 Chai\\u00A0: L\\u2019indulgence ultime pour vos sens\",\n    \"Th\u00e9 cha\u00ef\\u00A0: Le plaisir ultime pour tous vos sens\",\n  ],\n  [\n    \"Th\u00e9 Chai\\u00A0: Une \u00e9vasion douce du quotidien\",\n    \"Th\u00e9 cha\u00ef\\u00A0: \u00c9vadez-vous en douceur\",\n  ],\n  [\n    \"Th\u00e9 Chai\\u00A0: Partager la chaleur, partager l\\u2019amour\",\n    \"Th\u00e9 cha\u00ef\\u00A0: offrez de la chaleur, offrez de l\\u2019amour\",\n  ],\n  [\n    \"Th\u00e9 Chai\\u00A0: Traitez-vous \u00e0 quelque chose de sp\u00e9cial\",\n    \"Th\u00e9 cha\u00ef\\u00A0: Faites-vous plaisir\",\n  ],\n];\n\nconst replacementMap = new Map(replacements);\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  const newText = replacementMap.get(paragraph.text);\n  if (newText !== undefined) {\n    paragraph.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the French \"Th\u00e9 Chai\" slogans to the revised \"Th\u00e9 cha\u00ef\" copy.\n# Each slogan is a unique, single-run paragraph, so Find/Replace on the\n# full old sentence -> full new sentence (case-sensitive, no wildcards)\n# unambiguously targets the right run without touching its formatting.\n$d = $word.ActiveDocument\n$pairs = @(\n  @(\"Th\u00e9 Chai\u00a0: L\u2019\u00e9pices de la vie\", \"Th\u00e9 cha\u00ef\u00a0: L\u2019\u00e9pice de la vie\"),\n  @(\"Th\u00e9 Chai\u00a0: Un monde de saveur dans une tasse\", \"Th\u00e9 cha\u00ef\u00a0: Un monde de saveur dans une tasse\"),\n  @(\"Th\u00e9 Chai\u00a0: D\u00e9couvrir la magie de l\u2019Inde\", \"Th\u00e9 cha\u00ef\u00a0: D\u00e9couvrez la magie de l\u2019Inde\"),\n  @(\"Th\u00e9 Chai\u00a0: Le m\u00e9lange parfait de sant\u00e9 et de plaisir\", \"Th\u00e9 cha\u00ef\u00a0: L\u2019\u00e9quilibre parfait du bien-\u00eatre et du plaisir\"),\n  @(\"Th\u00e9 Chai\u00a0: Plus que le th\u00e9, un mode de vie\", \"Th\u00e9 cha\u00ef\u00a0: Plus qu\u2019un th\u00e9, un mode de vie\"),\n  @(\"Th\u00e9 Chai\u00a0: Une boisson pour toutes les saisons et les raisons\", \"Th\u00e9 cha\u00ef\u00a0: Une boisson pour toutes les saisons et toutes les raisons\"),\n  @(\"Th\u00e9 Chai\u00a0: L\u2019indulgence ultime pour vos sens\", \"Th\u00e9 cha\u00ef\u00a0: Le plaisir ultime pour tous vos sens\"),\n  @(\"Th\u00e9 Chai\u00a0: Une \u00e9vasion douce du quotidien\", \"Th\u00e9 cha\u00ef\u00a0: \u00c9vadez-vous en douceur\"),\n  @(\"Th\u00e9 Chai\u00a0: Partager la chaleur, partager l\u2019amour\", \"Th\u00e9 cha\u00ef\u00a0: offrez de la chaleur, offrez de l\u2019amour\"),\n  @(\"Th\u00e9 Chai\u00a0: Traitez-vous \u00e0 quelque chose de sp\u00e9cial\", \"Th\u00e9 cha\u00ef\u00a0: Faites-vous plaisir\"),\n)\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.Text = $pair[1]\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Execute([ref]$find.Text, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$find.Replacement.Text, [ref]2) | Out-Null\n}"}
